$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column I (row 3 to row 12) formatting into column J, then set the new values
$ws.Range("I3:I12").Copy() | Out-Null
$ws.Range("J3:J12").PasteSpecial() | Out-Null

$ws.Range("J4").Value = 2020
$ws.Range("J5").Value = 253.27664777870578
$ws.Range("J7").Value = 93.236077839070575
$ws.Range("J8").Value = 160
$ws.Range("J10").Value = 69
$ws.Range("J11").Value = 48.5
$ws.Range("J12").Value = 22.8

$ws.Range("J3").Select()
